$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = 44719
$ws.Range("M2").Value = 50
$ws.Range("N2").Value = 20000
$ws.Range("O2").Value = 21000
$ws.Range("P2").Value = 20400
$ws.Range("R2").Value = "Provincia de Limarí"
$ws.Range("S2").Value = 1133

# Row 3
$ws.Range("D3").Value = 45106
$ws.Range("M3").Value = 50
$ws.Range("N3").Value = 10000
$ws.Range("O3").Value = 10000
$ws.Range("P3").Value = 10000
$ws.Range("R3").Value = "Región de O'Higgins"
$ws.Range("S3").Value = 556

# Row 4
$ws.Range("D4").Value = 44320
$ws.Range("M4").Value = 50
$ws.Range("N4").Value = 18000
$ws.Range("O4").Value = 20000
$ws.Range("P4").Value = 18800
$ws.Range("R4").Value = "Provincia de Limarí"
$ws.Range("S4").Value = 1044

# Row 5
$ws.Range("D5").Value = 44362
$ws.Range("M5").Value = 100
$ws.Range("N5").Value = 19000
$ws.Range("O5").Value = 20000
$ws.Range("P5").Value = 19500
$ws.Range("R5").Value = "Provincia de Curicó"
$ws.Range("S5").Value = 1083

# Row 6
$ws.Range("D6").Value = 45084
$ws.Range("M6").Value = 100
$ws.Range("N6").Value = 17000
$ws.Range("O6").Value = 18000
$ws.Range("P6").Value = 17500
$ws.Range("R6").Value = "Región de O'Higgins"
$ws.Range("S6").Value = 972
